$p = $ppt.ActivePresentation
$nm = $p.NotesMaster
$tcs = $nm.Theme.ThemeColorScheme
$c1 = $tcs.Item(1)
$c1.RGB = 0x123456
Write-Output ("Set idx1 via NotesMaster to RGB=" + $c1.RGB)
